$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row before the existing row 137 (shifts old rows 137..172 down to 138..173,
# and the sheet dimension grows from A1:R172 to A1:R173 automatically).
$ws.Rows(137).Insert()

# Fill the newly-inserted row 137 with a new "Albahaca" price observation.
$ws.Range("A137").Value = 4
$ws.Range("B137").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C137").Value = "Los Lagos"
$ws.Range("D137").Value = Get-Date -Year 2023 -Month 2 -Day 28 -Hour 0 -Minute 0 -Second 0
$ws.Range("E137").Value = 10
$ws.Range("F137").Value = 100112052
$ws.Range("G137").Value = "Albahaca"
$ws.Range("H137").Value = "Sin especificar"
$ws.Range("I137").Value = "Primera"
$ws.Range("J137").Value = 90
$ws.Range("K137").Value = 5000
$ws.Range("L137").Value = 6000
$ws.Range("M137").Value = 5500
$ws.Range("N137").Value = "`$/docena de matas"
$ws.Range("O137").Value = "Región Metropolitana"
$ws.Range("P137").Value = 917
$ws.Range("Q137").Value = 6
$ws.Range("R137").Value = "Hortaliza"
